$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = 0.4270390154741289
$ws.Range("J19").Value = 0.2157505895368739
$ws.Range("K19").Value = -0.1957046463614244
$ws.Range("L19").Value = 2.961627167294484

$ws.Range("I20").Value = 0.7259974791323728
$ws.Range("J20").Value = 0.4796735828753367
$ws.Range("K20").Value = 0.2329828474207306
$ws.Range("L20").Value = 2.229085699281371
